# Update the "grid_cell" column (AG) on the "solar" sheet of the
# SubRES_REZoning_Sol_Win_and_Hydro workbook.
#
# The underlying SubRES generator re-derived the grid-cell assignment for
# the distributed solar/wind connection processes, so the "grid_cell"
# values shown in column AG (rows 4-25 of the "solar" sheet) need to be
# refreshed to the newly computed cell numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

# New grid_cell (column AG) values for rows 4-25, in row order.
$gridCells = @(
    "CHE_11",
    "CHE_15",
    "CHE_0",
    "CHE_3",
    "CHE_10",
    "CHE_22",
    "CHE_7",
    "CHE_20",
    "CHE_1",
    "CHE_24",
    "CHE_8",
    "CHE_5",
    "CHE_13",
    "CHE_17",
    "CHE_19",
    "CHE_23",
    "CHE_14",
    "CHE_18",
    "CHE_12",
    "CHE_9",
    "CHE_21",
    "CHE_4"
)

$startRow = 4
$col = 33  # column AG

for ($i = 0; $i -lt $gridCells.Length; $i++) {
    $ws.Cells.Item($startRow + $i, $col).Value = $gridCells[$i]
}
